$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new values look like plain numbers need to be forced to
# Text format first, so Excel stores them as literal strings (matching the source
# site's formatted-text export) instead of auto-converting them to numeric cells.
$textCells = @("D5", "D6", "D7", "D8", "D11", "D12", "D18", "D20", "D21", "D22", "D24", "D30", "D33", "D34", "D36", "D37", "D39", "D40", "D42", "D44", "D45", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.680.46'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '2.899.03'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '530.21'
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('D6').Value = '143.02'
$ws.Range('E6').Value = '  -6.24%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').Value = '2.903.06'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('E10').Value = '  -2.81%  '
$ws.Range('D11').Value = '5.96'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = '0.363'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').Value = '3.402.45'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = '60.635.23'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('D17').Value = '2.898.13'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').Value = '0.0000143'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = '11.69'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').Value = '364.06'
$ws.Range('E21').Value = '  -4.50%  '
$ws.Range('D22').Value = '6.62'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = '64.43'
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').Value = '3.019.39'
$ws.Range('E25').Value = '  -2.99%  '
$ws.Range('E26').Value = '  -2.85%  '
$ws.Range('E27').Value = '  -4.23%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0864'
$ws.Range('E29').Value = '  -6.94%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '7.78'
$ws.Range('E30').Value = '  -5.90%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').Value = '19.70'
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('D34').Value = '147.25'
$ws.Range('E34').Value = '  -7.03%  '
$ws.Range('D36').Value = '5.57'
$ws.Range('E36').Value = '  -6.83%  '
$ws.Range('D37').Value = '0.994'
$ws.Range('E37').Value = '  -6.27%  '
$ws.Range('E38').Value = '  -5.35%  '
$ws.Range('D39').Value = '37.70'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('D40').Value = '1.49'
$ws.Range('E40').Value = '  -5.11%  '
$ws.Range('D41').Value = '2.319.03'
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('D42').Value = '3.68'
$ws.Range('E42').Value = '  -4.71%  '
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D44').Value = '0.0581'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').Value = '20.57'
$ws.Range('E45').Value = '  -7.19%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('D49').Value = '0.0934'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').Value = '10.33'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '18.47'
$ws.Range('E51').Value = '  -5.87%  '

Write-Output "Applied crypto list update"